$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The publish date for book #7 (row 8) has now been verified, so the row is
# pulled out of its old spot and re-entered two rows further down (row 10),
# leaving row 9 blank in between.

# 1. Remove the old row 8 entirely (shifts rows below it up - none here).
$ws.Rows.Item(8).Delete()

# 2. Re-enter the same book data on row 10.
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "20 practice sets"
$ws.Range("C10").Value = "book for airman technical trades exam"
$ws.Range("D10").Value = "Arihant 'Expert Team'"
$ws.Range("E10").Value = 140
$ws.Range("F10").Value = 9780857503572
$ws.Range("G10").Value = "085750357X"
$ws.Range("H10").Value = "22\03\2004"
$ws.Range("I10").Value = "front.jpg"
$ws.Range("J10").Value = "back.jpg"
$ws.Range("K10").Value = 300
$ws.Range("L10").Value = "supporting.jpg"

# Match the workbook's existing cell formatting (9pt Calibri, style index 1).
$ws.Range("A10:L10").Font.Size = 9

# 3. Update the active selection to the newly (re)entered row.
$ws.Range("A10:L10").Select()
